$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 4-11 got cyclically rotated:
#   new row 4  <= old row 10
#   new row 5  <= old row 11
#   new row 6  <= old row 4
#   new row 7  <= old row 5
#   new row 8  <= old row 6
#   new row 9  <= old row 7
#   new row 10 <= old row 8
#   new row 11 <= old row 9
#
# Columns I, J, K, AF, AT, AY are empty placeholder cells in every one of
# these rows and never change, so they are left untouched entirely.
# Columns L, M, N, AC are empty-placeholder-or-content cells whose mere
# presence differs row to row, so after the bulk copy they get explicit
# add/remove/content fix-ups. Everything else with real content lives in
# contiguous blocks (A:H, P:W, Y:AB, AD:AE, AG:AG, AW:AX) that are safe to
# bulk-copy via .Value2 arrays.
#
# Y/AA hold literal "yyyy-mm-dd" text; Excel "helpfully" reinterprets such
# text as a real date when assigned via .Value2, so those two columns are
# re-applied afterwards as plain text.

function Get-RowBlocks($row) {
    return @{
        AH  = $ws.Range("A$row" + ":H$row").Value2
        PW  = $ws.Range("P$row" + ":W$row").Value2
        YAB = $ws.Range("Y$row" + ":AB$row").Value2
        ADE = $ws.Range("AD$row" + ":AE$row").Value2
        AG  = $ws.Range("AG$row").Value2
        AWX = $ws.Range("AW$row" + ":AX$row").Value2
    }
}

function Set-RowBlocks($row, $blocks) {
    $ws.Range("A$row" + ":H$row").Value2   = $blocks.AH
    $ws.Range("P$row" + ":W$row").Value2   = $blocks.PW
    $ws.Range("Y$row" + ":AB$row").Value2  = $blocks.YAB
    $ws.Range("AD$row" + ":AE$row").Value2 = $blocks.ADE
    $ws.Range("AG$row").Value2             = $blocks.AG
    $ws.Range("AW$row" + ":AX$row").Value2 = $blocks.AWX
}

# 1) Snapshot every source row's data (and its Startdatum/Slutdatum text)
#    before any writes happen.
$srcRows  = @(10, 11, 4, 5, 6, 7, 8, 9)
$dstRows  = @(4, 5, 6, 7, 8, 9, 10, 11)

$blocks    = @{}
$startDate = @{}
$endDate   = @{}
foreach ($r in $srcRows) {
    $blocks[$r]    = Get-RowBlocks $r
    $startDate[$r] = $ws.Range("Y$r").Value2
    $endDate[$r]   = $ws.Range("AA$r").Value2
}

# 2) Write the rotated rows.
for ($i = 0; $i -lt $dstRows.Length; $i++) {
    $dst = $dstRows[$i]
    $src = $srcRows[$i]
    Set-RowBlocks $dst $blocks[$src]
}

# 3) Restore Y/AA as literal text (undo any date reinterpretation).
for ($i = 0; $i -lt $dstRows.Length; $i++) {
    $dst = $dstRows[$i]
    $src = $srcRows[$i]

    $yCell = $ws.Range("Y$dst")
    $yCell.NumberFormat = "@"
    $yCell.Value2 = $startDate[$src]
    $yCell.ClearFormats()

    $aaCell = $ws.Range("AA$dst")
    $aaCell.NumberFormat = "@"
    $aaCell.Value2 = $endDate[$src]
    $aaCell.ClearFormats()
}

# 4) Fix up the L/M/N/AC presence-only columns for the new row layout.
#    "present" rows are restored to an empty-but-present cell by copying
#    an already-empty cell over them (keeps cell alive with no value);
#    "absent" rows are cleared outright so the cell disappears entirely.
$present = @{
    "L"  = @(6, 8, 9, 10)
    "M"  = @(10)
    "N"  = @(6, 8, 9, 10)
    "AC" = @()
}
$allRows = @(4, 5, 6, 7, 8, 9, 10, 11)

foreach ($col in @("L", "M", "N", "AC")) {
    $keepRows = $present[$col]
    foreach ($r in $allRows) {
        $cell = $ws.Range("$col$r")
        if ($keepRows -contains $r) {
            $ws.Range("I$r").Copy($cell)
        } else {
            $cell.ClearContents()
        }
    }
}

# AC11 carries real text content ("moved" there from the old AC9).
$ws.Range("AC11").Value2 = "På björk och al"
